$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7-10 (old Sales Tax Details, Broker Interest Received, Bond Interest Paid, Bond Interest Received, Total)
$ws.Rows("7:10").Delete()

# Update row 2: Trading Profit (value unchanged)
$ws.Range("A2").Value = "Trading Profit"

# Update row 3: Net Dividend Income
$ws.Range("A3").Value = "Net Dividend Income"
$ws.Range("B3").Value = 4773.629295799999

# Update row 4: Net Interest Profit
$ws.Range("A4").Value = "Net Interest Profit"
$ws.Range("B4").Value = 13869.468039

# Update row 5: OpEx
$ws.Range("A5").Value = "OpEx"
$ws.Range("B5").Value = -147.444

# Update row 6: Total
$ws.Range("A6").Value = "Total"
$ws.Range("B6").Value = 26526.53446901
